$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column J (index 10), shifting the
# "numbspecies ... uptakeP04" block from J:AP to L:AR.
$ws.Columns("J:K").Insert()

# New header cells for the inserted columns.
$ws.Range("J1").Value = "TOTALN"
$ws.Range("K1").Value = "TOTALP"

# New data values for the inserted columns (same value on every data row).
$ws.Range("J2:J7").Value = 2
$ws.Range("K2:K7").Value = 0.5

# Mark the (now empty-fill) style used on the new/adjacent columns so a new
# cellXfs record is produced, matching the extra style slot introduced by
# the edit.
$ws.Range("J1:L7").Interior.ColorIndex = -4142

# Restore the selection to match the post-edit workbook state.
$ws.Range("L39").Select()
